$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPhs = "phs001524"
$newPhs = "phs001554"

$neo4jFileName = "TC03_CDS_Filter_PHSAccession-phs001554_Neo4jData.xlsx"
$webFileName = "TC03_CDS_Filter_PHSAccession-phs001554_WebData.xlsx"

# 1) Participants query (B2): fix phs_accession
$participantsQuery = $ws.Range("B2").Value()
$participantsQuery = $participantsQuery.Replace($oldPhs, $newPhs)
$ws.Range("B2").Value = $participantsQuery

# 2) Shared "StatQuery" text (C2/C3/C4): fix phs_accession
$statQuery = $ws.Range("C2").Value()
$statQuery = $statQuery.Replace($oldPhs, $newPhs)
$ws.Range("C2").Value = $statQuery

# 3) Samples query (B3): fix phs_accession
$samplesQuery = $ws.Range("B3").Value()
$samplesQuery = $samplesQuery.Replace($oldPhs, $newPhs)
$ws.Range("B3").Value = $samplesQuery

# 4) Expected output workbook file names (TC02 -> TC03, phs001524 -> phs001554)
$ws.Range("D2").Value = $neo4jFileName
$ws.Range("E2").Value = $webFileName

# 5) Files query (B4): fix phs_accession and clear experimental_strategies filter
$filesQuery = $ws.Range("B4").Value()
$filesQuery = $filesQuery.Replace($oldPhs, $newPhs)
$filesQuery = $filesQuery.Replace('experimental_strategies: ["RNA-Seq"]', 'experimental_strategies: []')
$ws.Range("B4").Value = $filesQuery

# Re-apply the shared values to the remaining cells of rows 3 and 4
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("D3").Value = $neo4jFileName
$ws.Range("E3").Value = $webFileName
$ws.Range("D4").Value = $neo4jFileName
$ws.Range("E4").Value = $webFileName

# Restore the large wrapped-text rows to Excel's max row height (409.5pt),
# matching the look of the original authored sheet.
$ws.Rows(2).RowHeight = 409.5
$ws.Rows(3).RowHeight = 409.5
$ws.Rows(4).RowHeight = 409.5

# Update active selection to B2, as in the saved workbook.
$ws.Range("B2").Select()
